# Update the SOHP MODS template: the <mods:dateCreated> element used to be
# written with encoding="iso8601"; switch it to encoding="w3cdtf" (cell S1
# on Sheet1 holds the opening tag text for that element).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").Value = '<mods:originInfo><mods:dateCreated encoding="w3cdtf">'

# Mirror the author's UI interaction: they had scrolled right (so column H
# is the left-most visible column) and the edited cell (S1) ended up
# selected/active.
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("S1").Select()
